# Auto-generated: apply cell-value updates per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.536.15"
$ws.Range("E2").Value = "  +2.27%  "
$ws.Range("D3").Value = "2.588.93"
$ws.Range("E3").Value = "  +0.20%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'507.14"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.72%  "
$ws.Range("D6").Value = "'153.98"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.12%  "
$ws.Range("E7").Value = "  +0.53%  "
$ws.Range("E8").Value = "  -7.34%  "
$ws.Range("D9").Value = "2.595.95"
$ws.Range("E9").Value = "  +0.54%  "
$ws.Range("E10").Value = "  +6.94%  "
$ws.Range("E11").Value = "  +1.00%  "
$ws.Range("D12").Value = "'0.348"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +2.47%  "
$ws.Range("E13").Value = "  +0.85%  "
$ws.Range("D14").Value = "3.041.22"
$ws.Range("E14").Value = "  +2.19%  "
$ws.Range("D15").Value = "60.487.28"
$ws.Range("E15").Value = "  +2.26%  "
$ws.Range("D16").Value = "'21.52"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.83%  "
$ws.Range("E17").Value = "  +2.16%  "
$ws.Range("D18").Value = "2.589.28"
$ws.Range("E18").Value = "  +0.74%  "
$ws.Range("E19").Value = "  +0.91%  "
$ws.Range("D20").Value = "'345.78"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +3.90%  "
$ws.Range("D21").Value = "'10.45"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.17%  "
$ws.Range("E22").Value = "  +1.57%  "
$ws.Range("D23").Value = "'0.997"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.71%  "
$ws.Range("D24").Value = "'60.01"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.58%  "
$ws.Range("E25").Value = "  +1.22%  "
$ws.Range("D26").Value = "'0.167"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.35%  "
$ws.Range("D27").Value = "'0.999"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.53%  "
$ws.Range("D28").Value = "0.0₃0847"
$ws.Range("E28").Value = "  +2.89%  "
$ws.Range("D29").Value = "'7.34"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.71%  "
$ws.Range("E30").Value = "  +0.30%  "
$ws.Range("D31").Value = "'19.35"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.28%  "
$ws.Range("D32").Value = "'153.72"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.45%  "
$ws.Range("D33").Value = "'1.56"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.53%  "
$ws.Range("E34").Value = "  +3.68%  "
$ws.Range("E35").Value = "  +2.76%  "
$ws.Range("E36").Value = "  +0.08%  "
$ws.Range("D37").Value = "'0.861"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +13.00%  "
$ws.Range("D38").Value = "'0.850"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.00%  "
$ws.Range("E39").Value = "  +0.60%  "
$ws.Range("E40").Value = "  +1.53%  "
$ws.Range("E41").Value = "  +2.04%  "
$ws.Range("D42").Value = "'297.12"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +2.42%  "
$ws.Range("E43").Value = "  -1.85%  "
$ws.Range("D44").Value = "'0.615"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.67%  "
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").Value = "'0.997"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.14%  "
$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").Value = "'0.0557"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.14%  "
$ws.Range("E47").Value = "  +3.47%  "
$ws.Range("E48").Value = "  +0.19%  "
$ws.Range("E49").Value = "  -0.95%  "
$ws.Range("E50").Value = "  +0.66%  "
$ws.Range("D51").Value = "2.003.61"
$ws.Range("E51").Value = "  +0.39%  "
